# Update cryptos list values (price and volume(1h) columns), and
# swap two row pairs (OKB/Hedera at rows 35-36, Stacks/WEMIXToken at rows 47-48)
# to reflect the latest data pulled from coinranking.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.456.60'
$ws.Range("E2").Value = '  -2.60%  '
$ws.Range("D3").Value = '3.337.86'
$ws.Range("E3").Value = '  -4.60%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '185.01'
$ws.Range("E5").Value = '  -7.03%  '
$ws.Range("D6").Value = '535.55'
$ws.Range("E6").Value = '  -2.08%  '
$ws.Range("D7").Value = '0.609'
$ws.Range("E7").Value = '  +1.07%  '
$ws.Range("D8").Value = '3.335.01'
$ws.Range("E8").Value = '  -4.13%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = '0.626'
$ws.Range("E10").Value = '  -3.81%  '
$ws.Range("D11").Value = '60.53'
$ws.Range("E11").Value = '  -4.51%  '
$ws.Range("E12").Value = '  -4.26%  '
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").Value = '9.22'
$ws.Range("E14").Value = '  -5.35%  '
$ws.Range("D15").Value = '3.867.87'
$ws.Range("E15").Value = '  -4.50%  '
$ws.Range("D16").Value = '3.340.18'
$ws.Range("E16").Value = '  -4.40%  '
$ws.Range("E17").Value = '  -4.20%  '
$ws.Range("E18").Value = '  -2.64%  '
$ws.Range("D19").Value = '65.242.43'
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").Value = '11.29'
$ws.Range("E20").Value = '  -3.88%  '
$ws.Range("E21").Value = '  -5.33%  '
$ws.Range("D22").Value = '378.64'
$ws.Range("E22").Value = '  -2.98%  '
$ws.Range("D23").Value = '3.87'
$ws.Range("E23").Value = '  -3.19%  '
$ws.Range("D24").Value = '11.42'
$ws.Range("E24").Value = '  -2.44%  '
$ws.Range("D25").Value = '81.73'
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").Value = '3.92'
$ws.Range("E26").Value = '  +5.30%  '
$ws.Range("D27").Value = '6.00'
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("D29").Value = '11.73'
$ws.Range("E29").Value = '  -3.27%  '
$ws.Range("D30").Value = '8.56'
$ws.Range("E30").Value = '  -1.97%  '
$ws.Range("D31").Value = '29.29'
$ws.Range("E31").Value = '  -4.78%  '
$ws.Range("D32").Value = '652.62'
$ws.Range("E32").Value = '  -3.65%  '
$ws.Range("E33").Value = '  -2.67%  '
$ws.Range("E34").Value = '  -2.05%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.107'
$ws.Range("E35").Value = '  -2.89%  '
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '59.91'
$ws.Range("E36").Value = '  -5.26%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("E38").Value = '  -0.36%  '
$ws.Range("D39").Value = '37.07'
$ws.Range("E39").Value = '  -4.19%  '
$ws.Range("D40").Value = '0.0₃0734'
$ws.Range("E40").Value = '  +9.12%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("D43").Value = '2.920.90'
$ws.Range("E43").Value = '  -4.87%  '
$ws.Range("D44").Value = '2.55'
$ws.Range("E44").Value = '  +1.85%  '
$ws.Range("D45").Value = '2.74'
$ws.Range("E45").Value = '  -8.68%  '
$ws.Range("D46").Value = '0.0406'
$ws.Range("E46").Value = '  +2.33%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '2.68'
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.92'
$ws.Range("E48").Value = '  +12.13%  '
$ws.Range("D49").Value = '2.66'
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("D51").Value = '3.02'
$ws.Range("E51").Value = '  +4.29%  '
